# New Test data is added
# Adds two new worksheets (Contactus, ChangePassword) with sample data,
# mirroring the existing Login/StudentForm/Signup "RunMode" style sheets,
# updates the Login sheet selection, and moves the active tab/selection
# to the newly-added ChangePassword sheet.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Update the Login sheet selection to a plain range selection.
# ----------------------------------------------------------------------
$login = $wb.Worksheets.Item("Login")
$login.Range("A1:A4").Select()

# ----------------------------------------------------------------------
# 2) Add the "Contactus" worksheet at the end of the workbook.
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$contact = $wb.Worksheets.Add($null, $lastSheet)
$contact.Name = "Contactus"

$contact.Range("A1").Value = 0
$contact.Range("B1").Value = 1
$contact.Range("C1").Value = 2
$contact.Range("D1").Value = 3
$contact.Range("E1").Value = 4
$contact.Range("F1").Value = 5
$contact.Range("G1").Value = 6
$contact.Range("I1").Value = 7
$contact.Range("J1").Value = 8

$contact.Range("A2").Value = "Contactus"

$contact.Range("A3").Value = "RunMode"
$contact.Range("B3").Value = "Fname"
$contact.Range("C3").Value = "lname"
$contact.Range("D3").Value = "mailid"
$contact.Range("E3").Value = "Ph"
$contact.Range("F3").Value = "Company"
$contact.Range("G3").Value = "City"
$contact.Range("H3").Value = "Zip"
$contact.Range("I3").Value = "Comments"

$contact.Range("A4").Value = "Contactus"
$contact.Range("B4").Value = "Manikanta"
$contact.Range("C4").Value = "Thota"
$contact.Range("D4").Value = "mani6747@gmail.com"
$contact.Range("E4").Value = "966-659-7666"
$contact.Range("F4").Value = "eBiz Solutions"
$contact.Range("G4").Value = "Memphis"
$contact.Range("H4").Value = 38119
$contact.Range("I4").Value = "This sample checking for contact us"

# Hyperlink + "Hyperlink" cell style for the e-mail address, matching the
# style already used for the same address on the Login sheet.
$contact.Hyperlinks.Add($contact.Range("D4"), "mailto:mani6747@gmail.com") | Out-Null
$login.Range("B4").Copy()
$contact.Range("D4").PasteSpecial(-4122) | Out-Null

$contact.Columns.Item(1).ColumnWidth = 11.86
$contact.Columns.Item(4).ColumnWidth = 20.71
$contact.Columns.Item(5).ColumnWidth = 12.43
$contact.Columns.Item(6).ColumnWidth = 13.57
$contact.Columns.Item(9).ColumnWidth = 32.57

$contact.Range("H7").Select()

# ----------------------------------------------------------------------
# 3) Add the "ChangePassword" worksheet at the end of the workbook.
# ----------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$changePw = $wb.Worksheets.Add($null, $lastSheet2)
$changePw.Name = "ChangePassword"

$changePw.Range("A1").Value = 0
$changePw.Range("B1").Value = 1
$changePw.Range("C1").Value = 2

$changePw.Range("A2").Value = "ChangePassword"

$changePw.Range("A3").Value = "RunMode"
$changePw.Range("B3").Value = "CurrentPassword"
$changePw.Range("C3").Value = "NewPassword"

$changePw.Range("A4").Value = "ChangePassword"
$changePw.Range("B4").Value = 1234567
$changePw.Range("C4").Value = "mani123"

$changePw.Columns.Item(1).ColumnWidth = 16.14
$changePw.Columns.Item(2).ColumnWidth = 16.71
$changePw.Columns.Item(3).ColumnWidth = 14.14
$changePw.Columns.Item(4).ColumnWidth = 16.71

# This is the last-touched / right-most sheet, so make it the active tab
# (also restores a plain "B7" selection like the source workbook).
$changePw.Range("B7").Select()
$changePw.Activate()
